# clean up electrofishing parser, #749
#
# The "Destination Tank:" input box (cols A:D of row 2) is removed from both
# sheets, and a new "Destination Pond" column is inserted right before the
# "Location Name" column (i.e. a new column H, pushing the old H onward one
# column to the right). The "Fishing seconds" header is re-cased to
# "Fishing Seconds", and the new Destination Pond column is populated on the
# sample-data sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Template"
$ws2 = $wb.Worksheets.Item(2)   # "Sample Data"

# ---- Sample Data sheet (row 4/5 hold real data, so do this one first) ----

# Insert a new column before the old "Location Name" column (H), shifting
# everything from H onward one column to the right.
$ws2.Columns.Item(8).Insert()

# Remove the old "Destination Tank:" label + its (now orphaned) input cells.
$ws2.Range("A2:D2").Clear()

# Re-case the "Fishing seconds" header.
$ws2.Range("S3").Value = "Fishing Seconds"

# New "Destination Pond" column header + data.
$ws2.Range("H3").Value = "Desitination Pond"
$ws2.Range("H4").Value = "LP21"
$ws2.Range("H5").Value = "LP22A"

# Column got inserted with no width metadata copied from the donor column;
# match the sample sheet's best-fit width for the new text.
$ws2.Columns.Item(8).ColumnWidth = 16.1666666

# ---- Template sheet ----

$ws1.Columns.Item(8).Insert()

$ws1.Range("A2:D2").Clear()

$ws1.Range("S3").Value = "Fishing Seconds"

$ws1.Range("H3").Value = "Destination Pond"

# ---- Selections (Sample Data first so Template ends up the active tab) ----

$ws2.Range("H6").Select()
$ws1.Range("C1").Select()
